# Add "Role" and "IsCampCommitee" columns to the staff list worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("staff")

# Header row
$ws.Range("D1").Value = "Role"
$ws.Range("E1").Value = "IsCampCommitee"

# Data rows
$ws.Range("D2").Value = "Student"
$ws.Range("E2").Value = $true

$ws.Range("D3").Value = "Staff"
$ws.Range("E3").Value = $false

$ws.Range("D4").Value = "Student"
$ws.Range("E4").Value = $false

$ws.Range("D5").Value = "Student"
$ws.Range("E5").Value = $false

$ws.Range("D6").Value = "Staff"
$ws.Range("E6").Value = $false

# Column E width: best-fit (auto-fit) the new column to its contents
$ws.Columns.Item(5).AutoFit()

# Update selection to match target (cursor moved to D7)
$ws.Range("D7").Select()
